$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts existing Tipo/single from D to E)
$ws.Columns("D").Insert()

# Set header for new column D ("MAE") matching the style of the other headers
$ws.Range("D1").Value = "MAE"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Borders.LineStyle = 1
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160

# Set the MAE value for row 2
$ws.Range("D2").Value = 0.9239883288152645
